# B6-PowerPoint.pptx edit replay
#
# 1) Three tables (on the slides holding the "Unit 3" revision-style grids)
#    get their table style switched from the custom "Table_0" style
#    ({FBAED18D-6D61-4C99-9FA6-7D85D17B9207}, defined in tableStyles.xml)
#    to the built-in PowerPoint table style {9E5E3E3B-69E4-47A6-9FB8-3CA95A3AAAD1}.
#
# 2) The deck's visible theme (ppt/theme/theme1.xml, "Integral" / "Red Violet")
#    is recoloured to the stock PowerPoint "Office" palette (the palette that,
#    before this edit, only lived in ppt/theme/theme2.xml, the Notes Master's
#    theme) - i.e. the 12 theme colours are swapped onto the main theme.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Table styles
# ---------------------------------------------------------------------------
$newTableStyle = "{9E5E3E3B-69E4-47A6-9FB8-3CA95A3AAAD1}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Theme colours -> swap the "Office" palette onto the main theme
# ---------------------------------------------------------------------------
# RGBColor.RGB is a BGR-packed Long (0xBBGGRR), not RGB - pack accordingly.
function ToBgr([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Office theme colours, in ThemeColorScheme.Colors(1..12) order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($k = 1; $k -le $officeColors.Count; $k++) {
    $tcs.Colors($k).RGB = ToBgr($officeColors[$k - 1])
}
